# Replace Product, IT, and Finance templates with correct industry-specific content
#
# This script edits several worksheets of the Resource/Staffing Plan workbook,
# swapping AI/ML-flavored wording for generic "Product" wording, and touches a
# handful of otherwise-empty rows (so they get materialized as bare <row r="N"/>
# elements, matching the target OOXML) by toggling a no-op row property.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Resource Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Resource Overview")

$ws1.Range("A2").Value = "Product Development Implementation Project"
$ws1.Range("B6").Value = "Enterprise Product Development Implementation"
$ws1.Range("A18").Value = "Data Science/Product"
$ws1.Range("G18").Value = "Product, Python, Statistics"

# Materialize empty row 13 (currently absent from the sparse XML).
$ws1.Rows.Item(13).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Detailed Staffing Plan"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Detailed Staffing Plan")

$ws2.Range("A1").Value = "DETProductLED STAFFING PLAN"

$ws2.Range("C9").Value = "Data Science/Product"
$ws2.Range("K9").Value = "Product, Deep Learning, Python"
$ws2.Range("P9").Value = "Product Lead"

$ws2.Range("C10").Value = "Data Science/Product"
$ws2.Range("K10").Value = "Product, Statistics, R/Python"

$ws2.Range("C11").Value = "Data Science/Product"
$ws2.Range("K11").Value = "Product, Python, Visualization"

$ws2.Range("B12").Value = "Product Engineer"
$ws2.Range("C12").Value = "Data Science/Product"
$ws2.Range("K12").Value = "ProductOps, Python, Cloud"

$ws2.Range("C13").Value = "Data Science/Product"

# Materialize empty row 2 (currently absent from the sparse XML).
$ws2.Rows.Item(2).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 3: "Resource Timeline"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Resource Timeline")

# Materialize empty rows 2 and 11.
$ws3.Rows.Item(2).OutlineLevel = 0
$ws3.Rows.Item(11).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 4: "Skills Matrix"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Skills Matrix")

$ws4.Range("D3").Value = "Product Innovation"

# Materialize empty rows 2 and 11.
$ws4.Rows.Item(2).OutlineLevel = 0
$ws4.Rows.Item(11).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 5: "Cost Analysis"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Cost Analysis")

$ws5.Range("A6").Value = "Data Science/Product"

# Materialize empty rows 2, 14 and 15.
$ws5.Rows.Item(2).OutlineLevel = 0
$ws5.Rows.Item(14).OutlineLevel = 0
$ws5.Rows.Item(15).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 6: "Resource Risk Assessment"
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Resource Risk Assessment")

$ws6.Range("B5").Value = "Team lacks required Product expertise"

# Materialize empty rows 2, 12 and 13.
$ws6.Rows.Item(2).OutlineLevel = 0
$ws6.Rows.Item(12).OutlineLevel = 0
$ws6.Rows.Item(13).OutlineLevel = 0
